# Update header labels (row 1)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "Static"
$ws.Range("C1").Value = "Adaptive"
$ws.Range("D1").Value = "Static - Adaptive"
$ws.Range("E1").Value = "Adaptive - Static"

# Fix the "Static - Adaptive" column (D) to hold the actual difference
# between the Static (B) and Adaptive (C) columns instead of the old
# "1-(CM_0/CM_3)" ratio values.
$ws.Range("D2").Value = 39.73999999999999
$ws.Range("D3").Value = 397.6800000000001
